$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Beispielfirma GmbH"
$ws.Range("B3").Value = "Bf GmbH"
$ws.Range("B4").Value = "Berufsgenossenschaft Nahrungsmittel"
$ws.Range("B5").Value = "BGN"

$ws.Range("B7").Select()
